$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 176; this shifts the existing rows 176-196
# down to 177-197 (and copies formatting, e.g. the date style in column D,
# from the row above).
$ws.Rows.Item(176).Insert()

# Populate the newly inserted row 176 with the new weekly price record.
# Columns that stay the same as the (old, now-shifted) row 176 - i.e.
# Mercado/Region/Codreg/Categoria/Variedad/Calidad/Unidad/Origen/Kg -
# are copied across; only the date + price/volume figures change.
$ws.Range("A176").Value = 8
$ws.Range("B176").Value = "Terminal La Palmera de La Serena"
$ws.Range("C176").Value = "Coquimbo"
$ws.Range("D176").Value = 44617
$ws.Range("E176").Value = 4
$ws.Range("F176").Value = 100112021
$ws.Range("G176").Value = "Ají"
$ws.Range("H176").Value = "Americana (o)"
$ws.Range("I176").Value = "Primera"
$ws.Range("J176").Value = 540
$ws.Range("K176").Value = 12000
$ws.Range("L176").Value = 13000
$ws.Range("M176").Value = 12500
$ws.Range("N176").Value = "$/caja 15 kilos"
$ws.Range("O176").Value = "Provincia de Limarí"
$ws.Range("P176").Value = 833
$ws.Range("Q176").Value = 15
$ws.Range("R176").Value = "Hortaliza"
